$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2. This shifts the existing row 2 -> row 3
# and the existing row 3 -> row 4 (matching the target layout).
$ws.Rows("2:2").Insert()

# The inserted row inherits formatting from the row above (bold header);
# reset it back to the plain "Normal" style used by the other data rows.
$ws.Range("A2:Z2").Style = "Normal"

# New row 2 data
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 65465
$ws.Range("C2").Value = "periodo5"
$ws.Range("D2").Value = 35435
$ws.Range("E2").Value = "Mexicali"
$ws.Range("F2").Value = 65465
$ws.Range("G2").Value = "FAC. DE CIENCIAS HUMANAS"
$ws.Range("H2").Value = 54654
$ws.Range("I2").Value = "Luis"
$ws.Range("J2").Value = "Lopez"
$ws.Range("K2").Value = "Doriga"
$ws.Range("L2").Value = 2
$ws.Range("M2").Value = "Masculino"
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = "Doctorado"
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 1
$ws.Range("S2").Value = "SAUZAL"
$ws.Range("T2").Value = "MEXICO"
$ws.Range("U2").Value = "BAJA CALIFORNIA"
$ws.Range("V2").Value = "espanol"
$ws.Range("W2").Value = 3
$ws.Range("X2").Value = 1
$ws.Range("Y2").Value = "Estancia de Investigacion"
$ws.Range("Z2").Value = "judith"

# Row 4 (originally the old row 3, now pushed down) gets entirely new data.
$ws.Range("A4").Value = 9
$ws.Range("B4").Value = 2022
# "1" must be stored as text (not a number) in this column, so force text
# with a leading apostrophe and then drop the resulting quote-prefix style
# so the cell format matches the rest of the plain data cells.
$ws.Range("C4").Value = "'1"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "Ensenada"
$ws.Range("F4").Value = 400
$ws.Range("G4").Value = "FACULTAD DE CIENCIAS"
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = "Juana "
$ws.Range("J4").Value = "de Arco"
$ws.Range("K4").Value = "salen"
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = "Femenino"
$ws.Range("N4").Value = 3
$ws.Range("O4").Value = "Maestria"
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = "Universidad de Oaxaca"
$ws.Range("T4").Value = "MEXICO"
$ws.Range("U4").Value = "Oaxaca"
$ws.Range("V4").Value = "ESPAÑOL"
$ws.Range("W4").Value = 1
$ws.Range("X4").Value = 1
$ws.Range("Y4").Value = "Docencia"
$ws.Range("Z4").Value = "jluna@uabc.edu.mx"
